# Update "想去人数" (interested-count) figures on the 展览 and 全部类型
# sheets to the freshly scraped numbers (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" -----------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 12650
$ws1.Range("F3").Value  = 606
$ws1.Range("F5").Value  = 17
$ws1.Range("F6").Value  = 278
$ws1.Range("F7").Value  = 396
$ws1.Range("F9").Value  = 12630
$ws1.Range("F10").Value = 20
$ws1.Range("F11").Value = 3129
$ws1.Range("F16").Value = 1196
$ws1.Range("F21").Value = 6113

# --- Sheet "全部类型" ---------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 12650
$ws4.Range("F3").Value  = 606
$ws4.Range("F5").Value  = 17
$ws4.Range("F6").Value  = 278
$ws4.Range("F8").Value  = 396
$ws4.Range("F10").Value = 12631
$ws4.Range("F11").Value = 20
$ws4.Range("F12").Value = 3129
$ws4.Range("F17").Value = 1196
$ws4.Range("F23").Value = 6113
